# Add new columns I (I0) and J (IF) to the sheet, matching the header/style
# pattern already used for columns B..H, and fill in values for rows 2-70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add "I0" in I1 and "IF" in J1 ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold font, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows 2-70: values for columns I (I0) and J (IF) ---
$ijValues = @{
    2 = @(6, 6)
    3 = @(8, 8)
    4 = @(6, 6)
    5 = @(7, 7)
    6 = @(7, 8)
    7 = @(8, 8)
    8 = @(7, 7)
    9 = @(8, 8)
    10 = @(8, 8)
    11 = @(7, 8)
    12 = @(8, 8)
    13 = @(9, 9)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(7, 8)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(7, 8)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(7, 7)
    24 = @(8, 8)
    25 = @(8, 8)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(7, 8)
    29 = @(7, 8)
    30 = @(7, 7)
    31 = @(10, 10)
    32 = @(7, 7)
    33 = @(8, 8)
    34 = @(6, 7)
    35 = @(9, 9)
    36 = @(7, 7)
    37 = @(8, 8)
    38 = @(8, 8)
    39 = @(8, 8)
    40 = @(8, 8)
    41 = @(8, 8)
    42 = @(6, 7)
    43 = @(6, 7)
    44 = @(7, 7)
    45 = @(10, 10)
    46 = @(9, 9)
    47 = @(8, 8)
    48 = @(7, 7)
    49 = @(9, 9)
    50 = @(7, 7)
    51 = @(6, 6)
    52 = @(7, 8)
    53 = @(8, 8)
    54 = @(9, 9)
    55 = @(6, 7)
    56 = @(8, 9)
    57 = @(6, 7)
    58 = @(6, 6)
    59 = @(8, 8)
    60 = @(7, 8)
    61 = @(8, 8)
    62 = @(9, 9)
    63 = @(6, 7)
    64 = @(9, 9)
    65 = @(5, 5)
    66 = @(6, 6)
    67 = @(8, 8)
    68 = @(6, 6)
    69 = @(5, 5)
    70 = @(3, 3)
}

for ($r = 2; $r -le 70; $r++) {
    $pair = $ijValues[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}

Write-Host "I0 and IF columns added"
